$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 278: 2021-05-16
Set-TextCell 278 1 "2021-05-16"
Set-TextCell 278 2 "overview"
Set-TextCell 278 3 "K02000001"
Set-TextCell 278 4 "United Kingdom"
$ws.Cells.Item(278, 5).Value = 4450777
$ws.Cells.Item(278, 6).Value = 1926
$ws.Cells.Item(278, 7).Value = 4
$ws.Cells.Item(278, 8).Value = 127679

# Row 279: 2021-05-17
Set-TextCell 279 1 "2021-05-17"
Set-TextCell 279 2 "overview"
Set-TextCell 279 3 "K02000001"
Set-TextCell 279 4 "United Kingdom"
$ws.Cells.Item(279, 5).Value = 4452756
$ws.Cells.Item(279, 6).Value = 1979
$ws.Cells.Item(279, 7).Value = 5
$ws.Cells.Item(279, 8).Value = 127684

# Row 280: 2021-05-18
Set-TextCell 280 1 "2021-05-18"
Set-TextCell 280 2 "overview"
Set-TextCell 280 3 "K02000001"
Set-TextCell 280 4 "United Kingdom"
$ws.Cells.Item(280, 5).Value = 4450392
$ws.Cells.Item(280, 6).Value = 2412
$ws.Cells.Item(280, 7).Value = 7
$ws.Cells.Item(280, 8).Value = 127691
